$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove all existing hyperlinks first (clean slate) so we can rebuild them
# in the correct final order/positions without leaving stale relationship entries.
$ws.Hyperlinks.Delete()

$data = @(
  @('2025-09-09 12:38:05', '【AI活用】社内業務改善・RAG検索システム構築の依頼', 'システム開発', '200,000 円 ~ 300,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5390045', 353, '🔥AI,Ai ◇業務改善'),
  @('2025-09-09 12:38:05', '初回 LINEとChatGPTをつなげて、顧客対応できるボットを作りたいです', 'システム開発', '50,000 円 ~ 100,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5390438', 323, '🔥GPT,ChatGPT'),
  @('2025-09-09 12:38:05', '【AIクローン制作】肉声・映像応答のプロデュース依頼', 'システム開発', '100,000 円 ~ 200,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5389313', 303, '🔥AI,Ai'),
  @('2025-09-09 12:38:05', '【せどり×ツール製作】APIを使用したせどりツールを製作できるエンジニアさんを募集します♪', 'システム開発', '20,000 円 ~ 50,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5217096', 243, '🔥API ◆ツール'),
  @('2025-09-09 12:38:05', '【急募】ショッピングアプリの自動購入ツール開発依頼', 'システム開発', '100,000 円 ~ 200,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5390401', 163, '◆ツール,開発 ◇アプリ'),
  @('2025-09-09 12:38:05', '【Power Apps開発】グローバルタレントマネジメントPoCアプリ構築依頼', 'システム開発', '100,000 円 ~ 200,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5390188', 93, '◆開発 ◇アプリ'),
  @('2025-09-09 12:38:05', '【急募】共同開発エンジニアメンバーを募集します!', 'システム開発', '50,000 円 ~ 100,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5390021', 68, '◆開発'),
  @('2025-09-09 12:38:05', '【仮想通貨】自動売買EAの開発依頼', 'システム開発', '20,000 円 ~ 50,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5389714', 63, '◆開発'),
  @('2025-09-09 12:38:05', 'お歳暮受注発注集計スプレッドシート開発依頼', 'システム開発', '5,000 円 ~ 10,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5389953', 60, '◆開発'),
  @('2025-09-09 12:38:05', '〖リモート可〗Delphiエンジニア募集', 'システム開発', '300,000 円 ~ 500,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5341051', 25, $null),
  @('2025-09-09 12:38:05', '限定公開 PR 限定公開の仕事', 'システム開発', '500,000 円 ~ 1,000,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5385681', 25, $null),
  @('2025-09-09 12:38:05', '【講師募集】Gensparkを使ったWEB構築チュートリアル募集', 'システム開発', '10,000 円 ~ 20,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5390165', 10, $null)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $fcell = $ws.Cells.Item($row, 6)
    $fcell.Value = $rec[5]
    $ws.Hyperlinks.Add($fcell, $rec[5]) | Out-Null
    $ws.Cells.Item($row, 7).Value = $rec[6]
    if ($null -ne $rec[7]) {
        $ws.Cells.Item($row, 8).Value = $rec[7]
    } else {
        $ws.Cells.Item($row, 8).Clear()
    }
}
